# Auto-generated Excel COM-interop script to apply scheduled-runner data refresh
# to the Ultros_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 295.05554
$ws.Range("I9").Value = 178.57143
$ws.Range("K9").Value = 178.57143
$ws.Range("M9").Value = -9.571429999999992
$ws.Range("H28").Value = 1227.5
$ws.Range("I28").Value = 1034.6875
$ws.Range("J28").Value = 1998.75
$ws.Range("K28").Value = 1034.6875
$ws.Range("L28").Value = 1998.75
$ws.Range("M28").Value = -549.6875
$ws.Range("N28").Value = -2968.75
$ws.Range("H40").Value = 12450
$ws.Range("I40").Value = 9900
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 9900
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -9725
$ws.Range("N40").Value = -15350
$ws.Range("H41").Value = 391.2143
$ws.Range("I41").Value = 211.54546
$ws.Range("J41").Value = 1050
$ws.Range("K41").Value = 211.54546
$ws.Range("L41").Value = 1050
$ws.Range("M41").Value = 228.45454
$ws.Range("N41").Value = -1930
$ws.Range("H43").Value = 4517.421
$ws.Range("J43").Value = 4107.706
$ws.Range("L43").Value = 4107.706
$ws.Range("N43").Value = -4245.706
$ws.Range("H62").Value = 9605.333000000001
$ws.Range("I62").Value = 6658
$ws.Range("J62").Value = 15500
$ws.Range("K62").Value = 6658
$ws.Range("L62").Value = 15500
$ws.Range("M62").Value = -6034
$ws.Range("N62").Value = -16748
$ws.Range("H65").Value = 9605.333000000001
$ws.Range("I65").Value = 6658
$ws.Range("J65").Value = 15500
$ws.Range("K65").Value = 33290
$ws.Range("L65").Value = 77500
$ws.Range("M65").Value = -30170
$ws.Range("N65").Value = -83740
$ws.Range("H75").Value = 64573.25
$ws.Range("J75").Value = 64573.25
$ws.Range("L75").Value = 64573.25
$ws.Range("N75").Value = -66445.25
$ws.Range("H78").Value = 64573.25
$ws.Range("J78").Value = 64573.25
$ws.Range("L78").Value = 193719.75
$ws.Range("N78").Value = -203079.75
$ws.Range("H86").Value = 5557.143
$ws.Range("I86").Value = 5750
$ws.Range("J86").Value = 5480
$ws.Range("K86").Value = 5750
$ws.Range("L86").Value = 5480
$ws.Range("M86").Value = -4627
$ws.Range("N86").Value = -7726
$ws.Range("H87").Value = 19833.334
$ws.Range("J87").Value = 19833.334
$ws.Range("L87").Value = 19833.334
$ws.Range("N87").Value = -22329.334
$ws.Range("H89").Value = 5557.143
$ws.Range("I89").Value = 5750
$ws.Range("J89").Value = 5480
$ws.Range("K89").Value = 28750
$ws.Range("L89").Value = 27400
$ws.Range("M89").Value = -23134
$ws.Range("N89").Value = -38632
$ws.Range("H90").Value = 19833.334
$ws.Range("J90").Value = 19833.334
$ws.Range("L90").Value = 59500.00199999999
$ws.Range("N90").Value = -71980.00199999999
$ws.Range("H92").Value = 299.2
$ws.Range("I92").Value = 252.4375
$ws.Range("J92").Value = 486.25
$ws.Range("K92").Value = 252.4375
$ws.Range("L92").Value = 486.25
$ws.Range("M92").Value = 995.5625
$ws.Range("N92").Value = -2982.25
$ws.Range("H95").Value = 53999.5
$ws.Range("J95").Value = 53999.5
$ws.Range("L95").Value = 53999.5
$ws.Range("N95").Value = -59491.5
$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("K98").Value = 1000
$ws.Range("M98").Value = 498
$ws.Range("H105").Value = 39833.75
$ws.Range("J105").Value = 39833.75
$ws.Range("L105").Value = 39833.75
$ws.Range("N105").Value = -46821.75
$ws.Range("H106").Value = 4248.2334
$ws.Range("I106").Value = 4749.952
$ws.Range("J106").Value = 3077.5557
$ws.Range("K106").Value = 4749.952
$ws.Range("L106").Value = 3077.5557
$ws.Range("M106").Value = -4118.952
$ws.Range("N106").Value = -4339.5557
$ws.Range("H107").Value = 749.4783
$ws.Range("I107").Value = 741.3333
$ws.Range("K107").Value = 741.3333
$ws.Range("M107").Value = 1178.6667
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""
$ws.Range("H109").Value = 29769.23
$ws.Range("J109").Value = 29769.23
$ws.Range("L109").Value = 29769.23
$ws.Range("N109").Value = -32543.23
$ws.Range("H111").Value = 1712.7333
$ws.Range("I111").Value = 1621.7
$ws.Range("J111").Value = 1894.8
$ws.Range("K111").Value = 4865.1
$ws.Range("L111").Value = 5684.4
$ws.Range("M111").Value = -1798.1
$ws.Range("N111").Value = -11818.4
$ws.Range("H114").Value = 37656
$ws.Range("I114").Value = 37656
$ws.Range("K114").Value = 37656
$ws.Range("M114").Value = -33317
$ws.Range("H116").Value = 6578.2915
$ws.Range("I116").Value = 4555.5557
$ws.Range("J116").Value = 7791.933
$ws.Range("K116").Value = 4555.5557
$ws.Range("L116").Value = 7791.933
$ws.Range("M116").Value = -1113.5557
$ws.Range("N116").Value = -14675.933
$ws.Range("H118").Value = 15697
$ws.Range("I118").Value = 15697
$ws.Range("K118").Value = 47091
$ws.Range("M118").Value = -45434
$ws.Range("H120").Value = 61500
$ws.Range("J120").Value = 61500
$ws.Range("L120").Value = 61500
$ws.Range("N120").Value = -71176
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H124").Value = 39833.332
$ws.Range("J124").Value = 39833.332
$ws.Range("L124").Value = 39833.332
$ws.Range("N124").Value = -49653.332
$ws.Range("H125").Value = 1696.25
$ws.Range("I125").Value = 1611.6666
$ws.Range("K125").Value = 14504.9994
$ws.Range("M125").Value = -12044.9994
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws.Range("H134").Value = 39666.668
$ws.Range("J134").Value = 39666.668
$ws.Range("L134").Value = 39666.668
$ws.Range("N134").Value = -49806.668
$ws.Range("H135").Value = 9914.375
$ws.Range("I135").Value = 4605
$ws.Range("K135").Value = 41445
$ws.Range("M135").Value = -38910
$ws.Range("H137").Value = 3808.5715
$ws.Range("I137").Value = 2316.0667
$ws.Range("J137").Value = 5530.6924
$ws.Range("K137").Value = 6948.2001
$ws.Range("L137").Value = 16592.0772
$ws.Range("M137").Value = -4398.2001
$ws.Range("N137").Value = -21692.0772
$ws.Range("H138").Value = 2089.6072
$ws.Range("I138").Value = 1484.5
$ws.Range("K138").Value = 4453.5
$ws.Range("M138").Value = 686.5
$ws.Range("H139").Value = 42627.145
$ws.Range("J139").Value = 42627.145
$ws.Range("L139").Value = 42627.145
$ws.Range("N139").Value = -52907.145
$ws.Range("H140").Value = 39769.23
$ws.Range("J140").Value = 39769.23
$ws.Range("L140").Value = 39769.23
$ws.Range("N140").Value = -50129.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 3783.3333
$ws.Range("I16").Value = 750
$ws.Range("J16").Value = 9850
$ws.Range("K16").Value = 750
$ws.Range("L16").Value = 9850
$ws.Range("M16").Value = -463
$ws.Range("N16").Value = -10424
$ws.Range("H19").Value = 1802
$ws.Range("I19").Value = 1103
$ws.Range("J19").Value = 3549.5
$ws.Range("K19").Value = 1103
$ws.Range("L19").Value = 3549.5
$ws.Range("M19").Value = -874
$ws.Range("N19").Value = -4007.5
$ws.Range("H32").Value = 3131.3438
$ws.Range("I32").Value = 3583.2
$ws.Range("J32").Value = 1517.5714
$ws.Range("K32").Value = 3583.2
$ws.Range("L32").Value = 1517.5714
$ws.Range("M32").Value = -3296.2
$ws.Range("N32").Value = -2091.5714
$ws.Range("H45").Value = 5869.2666
$ws.Range("I45").Value = 3306.6
$ws.Range("J45").Value = 10994.6
$ws.Range("K45").Value = 3306.6
$ws.Range("L45").Value = 10994.6
$ws.Range("M45").Value = -2929.6
$ws.Range("N45").Value = -11748.6
$ws.Range("H74").Value = 4132.5
$ws.Range("I74").Value = 3149.125
$ws.Range("K74").Value = 3149.125
$ws.Range("M74").Value = -2275.125
$ws.Range("H77").Value = 4132.5
$ws.Range("I77").Value = 3149.125
$ws.Range("K77").Value = 15745.625
$ws.Range("M77").Value = -11377.625
$ws.Range("H88").Value = 1522.5625
$ws.Range("J88").Value = 1422.1538
$ws.Range("L88").Value = 1422.1538
$ws.Range("N88").Value = -2234.1538
$ws.Range("H91").Value = 1522.5625
$ws.Range("J91").Value = 1422.1538
$ws.Range("L91").Value = 1422.1538
$ws.Range("N91").Value = -4230.1538
$ws.Range("H97").Value = 1371.3334
$ws.Range("I97").Value = 777.8148
$ws.Range("K97").Value = 777.8148
$ws.Range("M97").Value = -281.8148
$ws.Range("H98").Value = 36974.5
$ws.Range("J98").Value = 36974.5
$ws.Range("L98").Value = 36974.5
$ws.Range("N98").Value = -42964.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14999.667
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 19999.5
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 19999.5
$ws.Range("M20").Value = -4753
$ws.Range("N20").Value = -20493.5
$ws.Range("H76").Value = 31901.857
$ws.Range("J76").Value = 31901.857
$ws.Range("L76").Value = 31901.857
$ws.Range("N76").Value = -32531.857
$ws.Range("H79").Value = 31901.857
$ws.Range("J79").Value = 31901.857
$ws.Range("L79").Value = 31901.857
$ws.Range("N79").Value = -34085.857
$ws.Range("H86").Value = 19270894
$ws.Range("I86").Value = 41751370
$ws.Range("K86").Value = 41751370
$ws.Range("M86").Value = -41750247
$ws.Range("H89").Value = 19270894
$ws.Range("I89").Value = 41751370
$ws.Range("K89").Value = 208756850
$ws.Range("M89").Value = -208751234
$ws.Range("H105").Value = 1827.7391
$ws.Range("I105").Value = 1019.26666
$ws.Range("K105").Value = 1019.26666
$ws.Range("M105").Value = 727.73334
$ws.Range("H107").Value = 9046.9
$ws.Range("I107").Value = 8250
$ws.Range("J107").Value = 10906.333
$ws.Range("K107").Value = 8250
$ws.Range("L107").Value = 10906.333
$ws.Range("M107").Value = -6330
$ws.Range("N107").Value = -14746.333
$ws.Range("H125").Value = 39833.332
$ws.Range("J125").Value = 39833.332
$ws.Range("L125").Value = 39833.332
$ws.Range("N125").Value = -49673.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 245.70833
$ws.Range("I7").Value = 98.666664
$ws.Range("K7").Value = 98.666664
$ws.Range("M7").Value = 14.333336
$ws.Range("H8").Value = 900
$ws.Range("J8").Value = 900
$ws.Range("L8").Value = 900
$ws.Range("N8").Value = -1180
$ws.Range("H15").Value = 550
$ws.Range("I15").Value = 500
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 500
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = -330
$ws.Range("N15").Value = -940
$ws.Range("H16").Value = 4540.9
$ws.Range("I16").Value = 2844.1428
$ws.Range("J16").Value = 8500
$ws.Range("K16").Value = 2844.1428
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = -2557.1428
$ws.Range("N16").Value = -9074
$ws.Range("H19").Value = 1033.1818
$ws.Range("I19").Value = 1148.75
$ws.Range("J19").Value = 725
$ws.Range("K19").Value = 1148.75
$ws.Range("L19").Value = 725
$ws.Range("M19").Value = -978.75
$ws.Range("N19").Value = -1065
$ws.Range("H21").Value = 2206.5
$ws.Range("J21").Value = 4000
$ws.Range("L21").Value = 4000
$ws.Range("N21").Value = -4470
$ws.Range("H22").Value = 294
$ws.Range("J22").Value = 214.66667
$ws.Range("L22").Value = 214.66667
$ws.Range("N22").Value = -914.6666700000001
$ws.Range("H24").Value = 1033.1818
$ws.Range("I24").Value = 1148.75
$ws.Range("J24").Value = 725
$ws.Range("K24").Value = 1148.75
$ws.Range("L24").Value = 725
$ws.Range("M24").Value = -978.75
$ws.Range("N24").Value = -1065
$ws.Range("H31").Value = 2666.25
$ws.Range("I31").Value = 1776.1177
$ws.Range("K31").Value = 1776.1177
$ws.Range("M31").Value = -1481.1177
$ws.Range("H34").Value = 2666.25
$ws.Range("I34").Value = 1776.1177
$ws.Range("K34").Value = 1776.1177
$ws.Range("M34").Value = -1574.1177
$ws.Range("H58").Value = 2391.1155
$ws.Range("I58").Value = 1541.4615
$ws.Range("K58").Value = 1541.4615
$ws.Range("M58").Value = -1338.4615
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H74").Value = 65157
$ws.Range("J74").Value = 100314
$ws.Range("L74").Value = 100314
$ws.Range("N74").Value = -102062
$ws.Range("H77").Value = 65157
$ws.Range("J77").Value = 100314
$ws.Range("L77").Value = 300942
$ws.Range("N77").Value = -309678
$ws.Range("H86").Value = 18962.182
$ws.Range("I86").Value = 10289.8
$ws.Range("J86").Value = 26189.166
$ws.Range("K86").Value = 10289.8
$ws.Range("L86").Value = 26189.166
$ws.Range("M86").Value = -9166.799999999999
$ws.Range("N86").Value = -28435.166
$ws.Range("H89").Value = 18962.182
$ws.Range("I89").Value = 10289.8
$ws.Range("J89").Value = 26189.166
$ws.Range("K89").Value = 51449
$ws.Range("L89").Value = 130945.83
$ws.Range("M89").Value = -45833
$ws.Range("N89").Value = -142177.83
$ws.Range("H99").Value = 2700.7144
$ws.Range("I99").Value = 2479
$ws.Range("K99").Value = 2479
$ws.Range("M99").Value = -981
$ws.Range("H113").Value = 4540.9
$ws.Range("I113").Value = 2844.1428
$ws.Range("J113").Value = 8500
$ws.Range("K113").Value = 2844.1428
$ws.Range("L113").Value = 8500
$ws.Range("M113").Value = -674.1428000000001
$ws.Range("N113").Value = -12840
$ws.Range("H126").Value = 2700.7144
$ws.Range("I126").Value = 2479
$ws.Range("K126").Value = 7437
$ws.Range("M126").Value = -4967
$ws.Range("H132").Value = 3405
$ws.Range("I132").Value = 3355.4285
$ws.Range("K132").Value = 10066.2855
$ws.Range("M132").Value = -7536.2855
$ws.Range("H134").Value = 4147.75
$ws.Range("I134").Value = 4347.7144
$ws.Range("J134").Value = 2748
$ws.Range("K134").Value = 13043.1432
$ws.Range("L134").Value = 8244
$ws.Range("M134").Value = -10508.1432
$ws.Range("N134").Value = -13314
$ws.Range("H136").Value = 2391.1155
$ws.Range("I136").Value = 1541.4615
$ws.Range("K136").Value = 4624.3845
$ws.Range("M136").Value = -2074.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 829.94446
$ws.Range("J12").Value = 959.2308
$ws.Range("L12").Value = 2877.6924
$ws.Range("N12").Value = -3223.6924
$ws.Range("H36").Value = 2137.6667
$ws.Range("I36").Value = 2365.2
$ws.Range("K36").Value = 7095.599999999999
$ws.Range("M36").Value = -6926.599999999999
$ws.Range("H40").Value = 70
$ws.Range("J40").Value = 81.8
$ws.Range("L40").Value = 327.2
$ws.Range("N40").Value = -465.2
$ws.Range("H68").Value = 2406.625
$ws.Range("I68").Value = 1225
$ws.Range("J68").Value = 2800.5
$ws.Range("K68").Value = 3675
$ws.Range("L68").Value = 8401.5
$ws.Range("M68").Value = -2864
$ws.Range("N68").Value = -10023.5
$ws.Range("H69").Value = 4083.5
$ws.Range("I69").Value = 2333.6667
$ws.Range("J69").Value = 5833.3335
$ws.Range("K69").Value = 7001.000100000001
$ws.Range("L69").Value = 17500.0005
$ws.Range("M69").Value = -6190.000100000001
$ws.Range("N69").Value = -19122.0005
$ws.Range("H71").Value = 2406.625
$ws.Range("I71").Value = 1225
$ws.Range("J71").Value = 2800.5
$ws.Range("K71").Value = 11025
$ws.Range("L71").Value = 25204.5
$ws.Range("M71").Value = -6969
$ws.Range("N71").Value = -33316.5
$ws.Range("H72").Value = 4083.5
$ws.Range("I72").Value = 2333.6667
$ws.Range("J72").Value = 5833.3335
$ws.Range("K72").Value = 21003.0003
$ws.Range("L72").Value = 52500.0015
$ws.Range("M72").Value = -16947.0003
$ws.Range("N72").Value = -60612.0015
$ws.Range("H80").Value = 6824.9375
$ws.Range("J80").Value = 6719.9
$ws.Range("L80").Value = 20159.7
$ws.Range("N80").Value = -22031.7
$ws.Range("H83").Value = 6824.9375
$ws.Range("J83").Value = 6719.9
$ws.Range("L83").Value = 60479.1
$ws.Range("N83").Value = -69839.10000000001
$ws.Range("H87").Value = 499.5
$ws.Range("I87").Value = 499.5
$ws.Range("K87").Value = 1498.5
$ws.Range("M87").Value = -250.5
$ws.Range("H90").Value = 499.5
$ws.Range("I90").Value = 499.5
$ws.Range("K90").Value = 4495.5
$ws.Range("M90").Value = 1744.5
$ws.Range("H114").Value = 850.2857
$ws.Range("I114").Value = 658.6667
$ws.Range("K114").Value = 1976.0001
$ws.Range("M114").Value = 1277.9999
$ws.Range("H122").Value = 4188.2856
$ws.Range("J122").Value = 5387.1875
$ws.Range("L122").Value = 48484.6875
$ws.Range("N122").Value = -53384.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 22704.75
$ws.Range("I33").Value = 17000
$ws.Range("J33").Value = 24606.334
$ws.Range("K33").Value = 17000
$ws.Range("L33").Value = 24606.334
$ws.Range("M33").Value = -16748
$ws.Range("N33").Value = -25110.334
$ws.Range("H38").Value = 22997.334
$ws.Range("J38").Value = 22997.334
$ws.Range("L38").Value = 22997.334
$ws.Range("N38").Value = -23923.334
$ws.Range("H40").Value = 21000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""
$ws.Range("H44").Value = 17999.834
$ws.Range("J44").Value = 15999.667
$ws.Range("L44").Value = 15999.667
$ws.Range("N44").Value = -17191.667
$ws.Range("H52").Value = 31598.6
$ws.Range("J52").Value = 31598.6
$ws.Range("L52").Value = 31598.6
$ws.Range("N52").Value = -32116.6
$ws.Range("H55").Value = 11946.6
$ws.Range("J55").Value = 14866.5
$ws.Range("L55").Value = 14866.5
$ws.Range("N55").Value = -15520.5
$ws.Range("H63").Value = 24557
$ws.Range("J63").Value = 24557
$ws.Range("L63").Value = 24557
$ws.Range("N63").Value = -25929
$ws.Range("H66").Value = 24557
$ws.Range("J66").Value = 24557
$ws.Range("L66").Value = 73671
$ws.Range("N66").Value = -80535
$ws.Range("H80").Value = 117890.3
$ws.Range("I80").Value = 225622
$ws.Range("J80").Value = 10158.6
$ws.Range("K80").Value = 225622
$ws.Range("L80").Value = 10158.6
$ws.Range("M80").Value = -224624
$ws.Range("N80").Value = -12154.6
$ws.Range("H83").Value = 117890.3
$ws.Range("I83").Value = 225622
$ws.Range("J83").Value = 10158.6
$ws.Range("K83").Value = 1128110
$ws.Range("L83").Value = 50793
$ws.Range("M83").Value = -1123118
$ws.Range("N83").Value = -60777
$ws.Range("H97").Value = 606.2
$ws.Range("I97").Value = 376.5
$ws.Range("J97").Value = 759.3333
$ws.Range("K97").Value = 376.5
$ws.Range("L97").Value = 759.3333
$ws.Range("M97").Value = 119.5
$ws.Range("N97").Value = -1751.3333
$ws.Range("H109").Value = 48500.2
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080
$ws.Range("H113").Value = 8129.385
$ws.Range("I113").Value = 2728.2856
$ws.Range("J113").Value = 14430.667
$ws.Range("K113").Value = 2728.2856
$ws.Range("L113").Value = 14430.667
$ws.Range("M113").Value = -558.2856000000002
$ws.Range("N113").Value = -18770.667
$ws.Range("H132").Value = 7835.7837
$ws.Range("I132").Value = 7285.769
$ws.Range("J132").Value = 9135.817999999999
$ws.Range("K132").Value = 21857.307
$ws.Range("L132").Value = 27407.454
$ws.Range("M132").Value = -19327.307
$ws.Range("N132").Value = -32467.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3469.0625
$ws.Range("I68").Value = 3161.111
$ws.Range("K68").Value = 3161.111
$ws.Range("M68").Value = -2412.111
$ws.Range("H71").Value = 3469.0625
$ws.Range("I71").Value = 3161.111
$ws.Range("K71").Value = 15805.555
$ws.Range("M71").Value = -12061.555
$ws.Range("H82").Value = 37038836
$ws.Range("I82").Value = 58825504
$ws.Range("J82").Value = 1495
$ws.Range("K82").Value = 58825504
$ws.Range("L82").Value = 1495
$ws.Range("M82").Value = -58825143
$ws.Range("N82").Value = -2217
$ws.Range("H85").Value = 37038836
$ws.Range("I85").Value = 58825504
$ws.Range("J85").Value = 1495
$ws.Range("K85").Value = 58825504
$ws.Range("L85").Value = 1495
$ws.Range("M85").Value = -58824256
$ws.Range("N85").Value = -3991
$ws.Range("H100").Value = 279527.75
$ws.Range("I100").Value = 279527.75
$ws.Range("K100").Value = 279527.75
$ws.Range("M100").Value = -278986.75
$ws.Range("H102").Value = 200561
$ws.Range("J102").Value = 200561
$ws.Range("L102").Value = 200561
$ws.Range("N102").Value = -207051
$ws.Range("H125").Value = 39833.332
$ws.Range("J125").Value = 39833.332
$ws.Range("L125").Value = 39833.332
$ws.Range("N125").Value = -49673.332
$ws.Range("H132").Value = 3003.5
$ws.Range("I132").Value = 2551.2
$ws.Range("K132").Value = 7653.599999999999
$ws.Range("M132").Value = -5123.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2000
$ws.Range("I9").Value = 2000
$ws.Range("K9").Value = 2000
$ws.Range("M9").Value = -1860
$ws.Range("H46").Value = 59106.75
$ws.Range("I46").Value = 57999
$ws.Range("J46").Value = 59476
$ws.Range("K46").Value = 57999
$ws.Range("L46").Value = 59476
$ws.Range("M46").Value = -57768
$ws.Range("N46").Value = -59938
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H100").Value = 570.4
$ws.Range("I100").Value = 610.8
$ws.Range("J100").Value = 530
$ws.Range("K100").Value = 1221.6
$ws.Range("L100").Value = 1060
$ws.Range("M100").Value = -680.5999999999999
$ws.Range("N100").Value = -2142
$ws.Range("H102").Value = 62221
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 62221
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 62221
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -68711
$ws.Range("H107").Value = 821.8461
$ws.Range("I107").Value = 835.5
$ws.Range("J107").Value = 776.3333
$ws.Range("K107").Value = 2506.5
$ws.Range("L107").Value = 2328.9999
$ws.Range("M107").Value = -586.5
$ws.Range("N107").Value = -6168.9999
$ws.Range("H126").Value = 4875
$ws.Range("I126").Value = 4833.3335
$ws.Range("K126").Value = 14500.0005
$ws.Range("M126").Value = -12030.0005
$ws.Range("H132").Value = 5001.5
$ws.Range("I132").Value = 3781.875
$ws.Range("J132").Value = 8253.833000000001
$ws.Range("K132").Value = 11345.625
$ws.Range("L132").Value = 24761.499
$ws.Range("M132").Value = -8815.625
$ws.Range("N132").Value = -29821.499
$ws.Range("H134").Value = 59106.75
$ws.Range("I134").Value = 57999
$ws.Range("J134").Value = 59476
$ws.Range("K134").Value = 173997
$ws.Range("L134").Value = 178428
$ws.Range("M134").Value = -171462
$ws.Range("N134").Value = -183498
$ws.Range("H136").Value = 66362.664
$ws.Range("I136").Value = 81474.164
$ws.Range("J136").Value = 5916.6665
$ws.Range("K136").Value = 244422.492
$ws.Range("L136").Value = 17749.9995
$ws.Range("M136").Value = -241872.492
$ws.Range("N136").Value = -22849.9995
